$wb = $excel.ActiveWorkbook

# --- Add the four new worksheets, in order, after "weights" ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetDefaultSkill = $wb.Worksheets.Add($null, $afterSheet)
$sheetDefaultSkill.Name = "defaultSkill"
$afterSheet = $sheetDefaultSkill
$sheetBarbarian = $wb.Worksheets.Add($null, $afterSheet)
$sheetBarbarian.Name = "barbarian"
$afterSheet = $sheetBarbarian
$sheetScout = $wb.Worksheets.Add($null, $afterSheet)
$sheetScout.Name = "scout"
$afterSheet = $sheetScout
$sheetKnight = $wb.Worksheets.Add($null, $afterSheet)
$sheetKnight.Name = "knight"

# --- Populate skill name / weight columns ---
$namesA = @(
  "FantasySkills.Academia",
  "FantasySkills.Acrobatics",
  "FantasySkills.Administration",
  "FantasySkills.AnimalHandling",
  "FantasySkills.Athletics",
  "FantasySkills.Craft",
  "FantasySkills.Engineering",
  "FantasySkills.Singing",
  "FantasySkills.Dancing",
  "FantasySkills.Sculpting",
  "FantasySkills.Music",
  "FantasySkills.Storytelling",
  "FantasySkills.Puppetry",
  "FantasySkills.Painting",
  "FantasySkills.Humanities",
  "FantasySkills.Intrusion",
  "FantasySkills.Investigation",
  "FantasySkills.Medicine",
  "FantasySkills.Meditation",
  "FantasySkills.Mysticism",
  "FantasySkills.Persuasion",
  "FantasySkills.Psychology",
  "FantasySkills.ScienceAlchemy",
  "FantasySkills.ScienceMathematics",
  "FantasySkills.ScienceAstronomy",
  "FantasySkills.ScienceNature",
  "FantasySkills.ScienceGeology",
  "FantasySkills.Stealth",
  "FantasySkills.Streetwise",
  "FantasySkills.Survival",
  "FantasySkills.Tactics",
  "FantasySkills.Trickery",
  "FantasySkills.Vehicle",
  "CombatSkills.Archery",
  "CombatSkills.BladedWeapon",
  "CombatSkills.BluntWeapon",
  "CombatSkills.Crossbows",
  "CombatSkills.Firearms",
  "CombatSkills.Shield",
  "CombatSkills.UnarmedStrikes",
  "CombatSKills.Wrestling"
)
$valsDefaultSkill = @(
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0,
  0.0
)
$valsBarbarian = @(
  -0.01,
  0.01,
  -0.01,
  0.01,
  0.02,
  0.0,
  0.0,
  0.01,
  0.01,
  0.0,
  0.0,
  0.01,
  0.0,
  0.0,
  0.0,
  0.005,
  0.0,
  0.005,
  0.0,
  0.008,
  0.0,
  0.0,
  -0.008,
  -0.008,
  -0.008,
  0.0,
  -0.008,
  0.01,
  0.01,
  0.03,
  0.01,
  0.01,
  0.0,
  0.008,
  0.03,
  0.03,
  0.0,
  0.0,
  0.01,
  0.008,
  0.008
)
$valsScout = @(
  -0.01,
  0.01,
  -0.01,
  0.01,
  0.02,
  0.005,
  0.0,
  0.01,
  0.01,
  0.0,
  0.01,
  0.01,
  0.0,
  0.0,
  0.0,
  0.01,
  0.01,
  0.005,
  0.0,
  0.008,
  0.0,
  0.0,
  -0.008,
  -0.008,
  0.0,
  0.0,
  0.0,
  0.02,
  0.01,
  0.03,
  0.01,
  0.01,
  0.0,
  0.03,
  0.01,
  0.008,
  0.01,
  0.01,
  -0.02,
  0.008,
  0.008
)
$valsKnight = @(
  0.0,
  -0.008,
  0.01,
  0.02,
  0.02,
  0.005,
  0.0,
  0.01,
  0.01,
  0.0,
  0.01,
  0.01,
  0.0,
  0.0,
  0.01,
  0.008,
  0.01,
  0.005,
  0.008,
  0.008,
  0.01,
  0.008,
  -0.008,
  -0.008,
  -0.008,
  -0.008,
  -0.008,
  0.0,
  0.0,
  0.0,
  0.02,
  0.0,
  0.01,
  0.008,
  0.04,
  0.04,
  0.008,
  0.008,
  0.015,
  0.01,
  0.01
)

for ($i = 0; $i -lt 41; $i++) {
  $r = $i + 1
  $sheetDefaultSkill.Cells.Item($r, 1).Value = $namesA[$i]
  $sheetDefaultSkill.Cells.Item($r, 2).Value = [double]$valsDefaultSkill[$i]
  $sheetBarbarian.Cells.Item($r, 1).Value = $namesA[$i]
  $sheetBarbarian.Cells.Item($r, 2).Value = [double]$valsBarbarian[$i]
  $sheetScout.Cells.Item($r, 1).Value = $namesA[$i]
  $sheetScout.Cells.Item($r, 2).Value = [double]$valsScout[$i]
  $sheetKnight.Cells.Item($r, 1).Value = $namesA[$i]
  $sheetKnight.Cells.Item($r, 2).Value = [double]$valsKnight[$i]
}

# --- Sum formulas (column D, row 1) on barbarian/scout/knight ---
$sheetBarbarian.Range("D1").Formula = "=SUM(B1:B41)"
$sheetScout.Range("D1").Formula = "=SUM(B1:B41)"
$sheetKnight.Range("D1").Formula = "=SUM(B1:B41)"

# --- Column widths (best effort given this runtime's width quantization) ---
$sheetDefaultSkill.Columns.Item(1).ColumnWidth = 32.0
$sheetBarbarian.Columns.Item(1).ColumnWidth = 30.833333333333332
$sheetBarbarian.Columns.Item(2).ColumnWidth = 12.666666666666666
$sheetBarbarian.Columns.Item(4).ColumnWidth = 17.666666666666668
$sheetScout.Columns.Item(1).ColumnWidth = 30.833333333333332
$sheetKnight.Columns.Item(1).ColumnWidth = 30.833333333333332

# --- View / selection state per sheet ---
$sheetDefaultSkill.Activate()
$sheetDefaultSkill.Range("A1:B41").Select() | Out-Null

$sheetBarbarian.Activate()
$sheetBarbarian.Range("A25").Select() | Out-Null

$sheetScout.Activate()
$sheetScout.Range("I35").Select() | Out-Null

$sheetKnight.Activate()
$sheetKnight.Range("G35").Select() | Out-Null

# --- weights sheet: select full row 6 ---
$wsWeights = $wb.Worksheets.Item("weights")
$wsWeights.Rows.Item(6).Select() | Out-Null

# --- Restore active tab to "weights" ---
$wsWeights.Activate()

